$d = $word.ActiveDocument

$d.Content.Find.Execute("[[PERSON_1]] – „s [[PERSON_1]]“, „o [[PERSON_2]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_1]] – „s [[PERSON_1]]“, „o [[PERSON_1]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_3]] – „k Evě Marečkové“, „u [[PERSON_4]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_2]] – „k [[PERSON_2]]“, „u [[PERSON_2]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_5]] – „pro [[PERSON_5]]“, „s [[PERSON_5]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_3]] – „pro [[PERSON_3]]“, „s [[PERSON_4]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_6]] – „s [[PERSON_6]]“, „o [[PERSON_7]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_5]] – „s [[PERSON_5]]“, „o [[PERSON_5]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_8]] – „u [[PERSON_9]]“, „k [[PERSON_10]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_6]] – „u [[PERSON_6]]“, „k [[PERSON_6]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_11]] – „o [[PERSON_11]]“, „se [[PERSON_11]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_7]] – „o [[PERSON_7]]“, „se [[PERSON_7]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_12]] – „k [[PERSON_13]]“, „u [[PERSON_14]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_8]] – „k [[PERSON_8]]“, „u [[PERSON_8]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_15]] – „s [[PERSON_15]]“, „o [[PERSON_16]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_9]] – „s [[PERSON_10]]“, „o [[PERSON_9]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_17]] – „u [[PERSON_18]]“, „s [[PERSON_17]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_11]] – „u [[PERSON_11]]“, „s [[PERSON_11]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_19]] – „s [[PERSON_19]]“, „k [[PERSON_19]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_12]] – „s [[PERSON_12]]“, „k [[PERSON_12]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_20]] – „s [[PERSON_20]]“, „o [[PERSON_21]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_13]] – „s [[PERSON_13]]“, „o [[PERSON_13]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_22]] – „ke [[PERSON_22]]“, „o [[PERSON_22]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_14]] – „ke [[PERSON_14]]“, „o [[PERSON_14]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_23]] – „o [[PERSON_24]]“, „s [[PERSON_23]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_15]] – „o [[PERSON_15]]“, „s [[PERSON_15]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_25]] – „u [[PERSON_26]]“, „s [[PERSON_25]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_16]] – „u [[PERSON_16]]“, „s [[PERSON_16]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_27]] – „ke [[PERSON_27]]“, „o [[PERSON_27]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_17]] – „ke [[PERSON_17]]“, „o [[PERSON_17]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_28]] – „s [[PERSON_28]]“, „o [[PERSON_28]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_18]] – „s [[PERSON_18]]“, „o [[PERSON_18]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_29]] – „s [[PERSON_29]]“, „o [[PERSON_29]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_19]] – „s [[PERSON_20]]“, „o [[PERSON_19]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_30]] – „k [[PERSON_30]]“, „od [[PERSON_31]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_21]] – „k [[PERSON_21]]“, „od [[PERSON_21]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_32]] – „o [[PERSON_32]]“, „s [[PERSON_32]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_22]] – „o [[PERSON_22]]“, „s [[PERSON_22]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_33]] – „o [[PERSON_33]]“, „se [[PERSON_33]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_23]] – „o [[PERSON_23]]“, „se [[PERSON_23]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_34]] – „s [[PERSON_34]]“, „u [[PERSON_35]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_24]] – „s [[PERSON_25]]“, „u [[PERSON_24]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_36]] – „o [[PERSON_37]]“, „s [[PERSON_36]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_26]] – „o [[PERSON_26]]“, „s [[PERSON_26]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_38]] – „k [[PERSON_39]]“, „o [[PERSON_39]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_27]] – „k [[PERSON_27]]“, „o [[PERSON_27]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_40]] – „se [[PERSON_40]]“, „o Soně Mikulkové“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_28]] – „se [[PERSON_28]]“, „o [[PERSON_28]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_41]] – „o [[PERSON_42]]“, „s [[PERSON_41]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_29]] – „o [[PERSON_29]]“, „s [[PERSON_29]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_43]] – „s [[PERSON_43]]“, „o [[PERSON_43]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_30]] – „s [[PERSON_30]]“, „o [[PERSON_30]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_44]] – „k [[PERSON_45]]“, „s [[PERSON_44]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_31]] – „k [[PERSON_31]]“, „s [[PERSON_31]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_47]] – „od [[PERSON_48]]“, „s [[PERSON_49]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_32]] – „s [[PERSON_32]]“, „o [[PERSON_32]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_50]] – „k [[PERSON_51]]“, „o [[PERSON_51]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_33]] – „od [[PERSON_33]]“, „s [[PERSON_34]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_52]] – „o [[PERSON_53]]“, „s [[PERSON_52]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_35]] – „k [[PERSON_35]]“, „o [[PERSON_35]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_54]] – „s [[PERSON_54]]“, „o [[PERSON_55]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_36]] – „o [[PERSON_36]]“, „s [[PERSON_36]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_56]] – „s [[PERSON_56]]“, „o [[PERSON_56]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_37]] – „s [[PERSON_37]]“, „o [[PERSON_37]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_57]] – „k [[PERSON_57]]“, „s [[PERSON_57]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_38]] – „s [[PERSON_38]]“, „o [[PERSON_38]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_58]] – „pro [[PERSON_59]]“, „o [[PERSON_60]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_39]] – „k [[PERSON_39]]“, „s [[PERSON_39]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_61]] – „k [[PERSON_61]]“, „o [[PERSON_61]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_40]] – „pro [[PERSON_41]]“, „o [[PERSON_42]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_62]] – „o [[PERSON_63]]“, „s [[PERSON_62]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_43]] – „k [[PERSON_43]]“, „o [[PERSON_43]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_64]] – „s [[PERSON_64]]“, „o [[PERSON_65]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_44]] – „o [[PERSON_44]]“, „s [[PERSON_44]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_66]] – „s [[PERSON_66]]“, „o [[PERSON_66]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_45]] – „s [[PERSON_45]]“, „o [[PERSON_45]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_67]] – „u [[PERSON_68]]“, „o [[PERSON_68]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_47]] – „u [[PERSON_47]]“, „o [[PERSON_47]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_69]] – „se [[PERSON_69]]“, „o [[PERSON_69]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_48]] – „se [[PERSON_48]]“, „o [[PERSON_48]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_70]] – „o [[PERSON_71]]“, „s [[PERSON_70]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_49]] – „o [[PERSON_50]]“, „s [[PERSON_50]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_72]] – „k [[PERSON_73]]“, „o [[PERSON_73]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_51]] – „k [[PERSON_51]]“, „o [[PERSON_51]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_74]] – „o [[PERSON_75]]“, „s [[PERSON_74]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_52]] – „o [[PERSON_52]]“, „s [[PERSON_52]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_76]] – „s [[PERSON_76]]“, „o [[PERSON_76]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_53]] – „s [[PERSON_53]]“, „o [[PERSON_53]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_77]] – „s [[PERSON_77]]“, „o [[PERSON_77]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_54]] – „s [[PERSON_54]]“, „o [[PERSON_54]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_78]] – „o [[PERSON_78]]“, „s [[PERSON_78]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_55]] – „o [[PERSON_55]]“, „s [[PERSON_55]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_79]] – „s [[PERSON_80]]“, „o [[PERSON_81]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_56]] – „s [[PERSON_57]]“, „o [[PERSON_56]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_82]] – „o [[PERSON_83]]“, „s [[PERSON_82]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_58]] – „o [[PERSON_58]]“, „s [[PERSON_58]]“", 2) | Out-Null
$d.Content.Find.Execute("[[PERSON_84]] – „s [[PERSON_85]]“, „o [[PERSON_86]]“", $true, $false, $false, $false, $false, $true, 1, $false, "[[PERSON_59]] – „s [[PERSON_60]]“, „o [[PERSON_59]]“", 2) | Out-Null

Write-Output "Replaced 49 paragraphs"
